$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision on which movie to show on Friday ended without a clear choice.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired. The committee has decided to show both movies.`n"
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision-making process concluded without a selection for Friday's movie, so no action will be taken in this instance.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" about the movie for Friday.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for both movies.`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has been recorded as `"no decision.`"`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie will be shown on Friday.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has ended without an agreement, resulting in no decision being made.`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday resulted in no selection.`n"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has resulted in no agreement among the committee members.`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C17").Value = "MSG: None`n`nMSG: The conversation ended without a plan about what movie to play on Friday, so no decision can be made.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday has been recorded as no decision.`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday has resulted in no specific choice being made.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no consensus was reached regarding the movie to be shown on Friday.`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded, and no specific movie was selected for Friday's screening.`n"
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday resulted in no agreement. The conversation did not lead to a definitive choice.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no_decision,`" indicating that no consensus was reached regarding a movie for Friday.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was selected for Friday’s screening.`n"
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been finalized.`n"
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights to both movies, `"Barbie`" and `"Oppenheimer,`" as they were both agreed upon for showing in their entirety.`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired successfully.`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision process resulted in no movie being selected for Friday.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision has been made to not acquire any movie for Friday at this time, as the committee did not reach a conclusion.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been successfully recorded.`n"
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision process has concluded without reaching an agreement on which movie to screen on Friday.`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for the movie `"Barbie.`"`n"
$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie selection for Friday.`n"
$ws.Range("C36").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision has been recorded for acquiring the rights to both movies.`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding which movie to show on Friday.`n"
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has resulted in no conclusive agreement.`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was selected for showing on Friday.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no_decision.`" The committee did not reach an agreement on which movie to show on Friday.`n"
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Range("C44").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired.`n"
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision about Friday's movie was not reached, so no action has been taken.`n"
$ws.Range("C47").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision resulted in no choice being made for Friday's movie.`n"
$ws.Range("C49").Value = "MSG: None`n`nMSG: It seems there was no decision made regarding the movie to be shown on Friday, so I will record that as the outcome.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been made.`n"
$ws.Range("C51").Value = "MSG: None`n`nMSG: The decision concluded with no movie being selected.`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C53").Value = "MSG: None`n`nMSG: The decision process has concluded without agreement on a movie for Friday, resulting in no decision being made.`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has not been made, resulting in no acquisition of rights for either movie.`n"
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been recorded successfully.`n"
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C57").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie will be acquired for Friday's showing.`n"
$ws.Range("C59").Value = "MSG: None`n`nMSG: The decision resulted in no agreement about which movie to show on Friday.`n"
$ws.Range("C60").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired successfully.`n"
$ws.Range("C61").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("C62").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been made, and the conversation ended without a clear choice.`n"
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision about the movie to be shown on Friday was not made, as the committee did not arrive at a conclusion.`n"
$ws.Range("C65").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has not been made.`n"
$ws.Range("C66").Value = "MSG: None`n`nMSG: The committee did not make a decision on which movie to show on Friday.`n"

$ws.Range("D12").Value = "Barbie_was_selected, "
$ws.Range("D17").Value = "no_decision, "
$ws.Range("D30").Value = "no_decision, "
$ws.Range("D34").Value = "Barbie_was_selected, "
$ws.Range("D42").Value = "no_decision, "
$ws.Range("D44").Value = "both_movies, "
